$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its text nature, matching the original
# inlineStr cells, so numeric-looking values (e.g. "582.32") are not
# auto-converted into real numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '67.196.23'
$ws.Range('E2').Value = '  +2.43%  '
$ws.Range('D3').Value = '3.106.13'
$ws.Range('E3').Value = '  +5.22%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '582.32'
$ws.Range('E5').Value = '  +2.63%  '
$ws.Range('D6').Value = '171.15'
$ws.Range('E6').Value = '  +7.72%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('D8').Value = '3.102.94'
$ws.Range('E8').Value = '  +5.16%  '
$ws.Range('D9').Value = '0.528'
$ws.Range('E9').Value = '  +2.15%  '
$ws.Range('D10').Value = '6.67'
$ws.Range('E10').Value = '  -0.96%  '
$ws.Range('E11').Value = '  +3.31%  '
$ws.Range('E12').Value = '  +5.46%  '
$ws.Range('D13').Value = '0.0000253'
$ws.Range('E13').Value = '  +3.23%  '
$ws.Range('D14').Value = '37.15'
$ws.Range('E14').Value = '  +9.11%  '
$ws.Range('E15').Value = '  -0.36%  '
$ws.Range('D16').Value = '3.624.61'
$ws.Range('E16').Value = '  +5.29%  '
$ws.Range('D17').Value = '67.183.53'
$ws.Range('E17').Value = '  +2.31%  '
$ws.Range('D18').Value = '7.26'
$ws.Range('E18').Value = '  +4.57%  '
$ws.Range('D19').Value = '3.106.83'
$ws.Range('E19').Value = '  +5.27%  '
$ws.Range('D20').Value = '16.14'
$ws.Range('E20').Value = '  +16.84%  '
$ws.Range('D21').Value = '473.82'
$ws.Range('E21').Value = '  +6.22%  '
$ws.Range('D22').Value = '0.720'
$ws.Range('E22').Value = '  +6.30%  '
$ws.Range('D23').Value = '7.57'
$ws.Range('E23').Value = '  +5.60%  '
$ws.Range('D24').Value = '84.06'
$ws.Range('E24').Value = '  +1.70%  '
$ws.Range('D25').Value = '2.38'
$ws.Range('E25').Value = '  +9.04%  '
$ws.Range('D26').Value = '13.00'
$ws.Range('E26').Value = '  +7.31%  '
$ws.Range('D27').Value = '10.34'
$ws.Range('E27').Value = '  +3.90%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('D29').Value = '8.16'
$ws.Range('E29').Value = '  +3.31%  '
$ws.Range('D30').Value = '2.45'
$ws.Range('E30').Value = '  +5.79%  '
$ws.Range('E31').Value = '  +4.85%  '
$ws.Range('D32').Value = '0.0000103'
$ws.Range('E32').Value = '  +5.13%  '
$ws.Range('D33').Value = '28.55'
$ws.Range('E33').Value = '  +4.54%  '
$ws.Range('E34').Value = '  +5.66%  '
$ws.Range('E35').Value = '  +0.12%  '
$ws.Range('E36').Value = '  +4.56%  '
$ws.Range('D37').Value = '5.96'
$ws.Range('E37').Value = '  +4.51%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = '2.13'
$ws.Range('E38').Value = '  +8.01%  '
$ws.Range('B39').Value = 'Arweave'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D39').Value = '47.36'
$ws.Range('E39').Value = '  +10.34%  '
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').Value = '0.322'
$ws.Range('E40').Value = '  +7.44%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').Value = '50.54'
$ws.Range('E41').Value = '  +2.93%  '
$ws.Range('E42').Value = '  +4.44%  '
$ws.Range('D43').Value = '2.92'
$ws.Range('E43').Value = '  +5.63%  '
$ws.Range('D44').Value = '8.75'
$ws.Range('E44').Value = '  +3.91%  '
$ws.Range('D45').Value = '396.99'
$ws.Range('E45').Value = '  +2.32%  '
$ws.Range('D46').Value = '0.0366'
$ws.Range('E46').Value = '  +3.67%  '
$ws.Range('D47').Value = '2.779.69'
$ws.Range('E47').Value = '  +1.87%  '
$ws.Range('D48').Value = '135.15'
$ws.Range('E48').Value = '  +3.09%  '
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('D50').Value = '24.90'
$ws.Range('E50').Value = '  +7.78%  '
$ws.Range('D51').Value = '2.26'
$ws.Range('E51').Value = '  +5.50%  '
